$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Přerozdělené body" values for the 12.týden column (I) for rows 4-7
$ws.Range("I4").Value = -3
$ws.Range("I5").Value = 2
$ws.Range("I6").Value = 2
$ws.Range("I7").Value = -1

# Update the "Přerozděleno" row (row 9): 10.týden (G9), 12.týden reason (H9) and value (I9)
$ws.Range("G9").Value = 4
$ws.Range("H9").Value = "-"
$ws.Range("I9").Value = 4

# Freeze panes at column B (split after column A), with selection in pane topRight
$ws.Range("B1").Select()
$excel.ActiveWindow.FreezePanes = $true
